$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$ws.Range('D2').Value = '29.611.51'
$ws.Range('E2').Value = '  +2.61%  '
$ws.Range('D3').Value = '1.862.57'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +2.73%  '
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.75%  '
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').Value = '1.861.32'
$ws.Range('E13').Value = '  +1.97%  '
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('E15').Value = '  +3.30%  '
$ws.Range('E16').Value = '  +2.66%  '
$ws.Range('D17').Value = '29.591.52'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').Value = '2.109.69'
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +3.33%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +2.93%  '
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E34').Value = '  +4.25%  '
$ws.Range('E35').Value = '  +5.02%  '
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').Value = '1.335.48'
$ws.Range('E38').Value = '  +10.84%  '
$ws.Range('E39').Value = '  +3.08%  '
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('E41').Value = '  +4.79%  '
$ws.Range('E42').Value = '  +14.54%  '
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E45').Value = '  +4.10%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E46').Value = '  +4.01%  '
$ws.Range('D47').Value = '2.010.25'
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('E49').Value = '  +4.13%  '
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('E51').Value = '  +2.30%  '

# Numeric-looking values that must stay as literal text (preserve exact
# formatting such as trailing zeros / leading zeros): force NumberFormat to
# text ("@") before the write, then restore the original style so no visible
# formatting/style change is introduced.
$textCells = @{
    'D4' = '0.9999'
    'D5' = '245.59'
    'D6' = '0.6998'
    'D8' = '0.07744'
    'D9' = '0.3071'
    'D10' = '23.64'
    'D11' = '0.07789'
    'D12' = '5.164'
    'D14' = '92.44'
    'D15' = '0.6934'
    'D16' = '6.603'
    'D18' = '0.000008359'
    'D20' = '242.51'
    'D21' = '12.78'
    'D23' = '7.627'
    'D25' = '0.1512'
    'D26' = '8.922'
    'D27' = '159.86'
    'D28' = '18.35'
    'D29' = '1.539'
    'D30' = '4.268'
    'D31' = '4.193'
    'D32' = '1.193'
    'D33' = '0.05104'
    'D34' = '0.7854'
    'D35' = '1.904'
    'D36' = '1.157'
    'D37' = '2.689'
    'D39' = '0.01880'
    'D40' = '2.733'
    'D41' = '0.9556'
    'D42' = '5.998'
    'D43' = '106.75'
    'D45' = '0.00000000127'
    'D46' = '9.786'
    'D48' = '0.5218'
    'D50' = '1.788'
    'D51' = '7.018'
}
foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.Style = $origStyle
}

Write-Host "Applied $($textCells.Count + 60) cell updates"